$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add a new row (row 6) with the same formatting as an existing full row
#    (row 2 has all 7 data columns populated) before touching values, so the
#    new row's cells inherit style "2" like the rest of the body rows.
# ---------------------------------------------------------------------------
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A6:G6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. The Copyright column (C) never had any content/style in rows 4 and 5 -
#    pull the formatting that already exists in C2:C3 down into them so the
#    new values get the same style ("s=2") as the rest of the column. Same
#    story for G5 (the Orientation column never had a cell there before).
# ---------------------------------------------------------------------------
$ws.Range("C2:C3").Copy() | Out-Null
$ws.Range("C4:C5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("G4").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Remove the old "Duration"/column H data - the table is now only A:G.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Filename"
$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Copyright"
$ws.Range("D1").Value = "Year"
$ws.Range("E1").Value = "Media"
$ws.Range("F1").Value = "Dimensions"
$ws.Range("G1").Value = "Orientation"

# ---------------------------------------------------------------------------
# 5. Row 2 - Hell Block
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Hell Block.jpg"
$ws.Range("B2").Value = "Hell Block"
$ws.Range("C2").Value = "Alyssa Taylor"
$ws.Range("D2").Value = 2013
$ws.Range("E2").Value = "Oil Paint on Canvas"
$ws.Range("F2").Value = "48""x36"""
$ws.Range("G2").Value = "Landscape"

# ---------------------------------------------------------------------------
# 6. Row 3 - Still Life in White #1 (new position)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Still Life in White #1.jpg"
$ws.Range("B3").Value = "Still Life in White #1"
$ws.Range("C3").Value = "Alyssa Taylor"
$ws.Range("D3").Value = 2013
$ws.Range("E3").Value = "Oil Paint on Canvas"
$ws.Range("F3").Value = "18""x24"""
$ws.Range("G3").Value = "Portrait"

# ---------------------------------------------------------------------------
# 7. Row 4 - Lucky #13
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Lucky #13.JPG"
$ws.Range("B4").Value = "Lucky #13"
$ws.Range("C4").Value = "Alyssa Taylor"
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = "Oil Paint on Canvas"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "Portrait"

# ---------------------------------------------------------------------------
# 8. Row 5 - Wisteria #1
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Wisteria #1.JPG"
$ws.Range("B5").Value = "Wisteria #1"
$ws.Range("C5").Value = "Alyssa Taylor"
$ws.Range("D5").Value = 2016
$ws.Range("E5").Value = "Synthetic Polymer on Canvas"
$ws.Range("F5").Clear() | Out-Null
$ws.Range("G5").Value = "Portrait"

# ---------------------------------------------------------------------------
# 9. Row 6 - Stella (new row)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Stella.jpg"
$ws.Range("B6").Value = "Stella"
$ws.Range("C6").Value = "Alyssa Taylor"
$ws.Range("D6").Value = 2016
$ws.Range("E6").Value = "Oil Paint on Canvas"
$ws.Range("F6").Clear() | Out-Null
$ws.Range("G6").Value = "Landscape"

# ---------------------------------------------------------------------------
# 10. Row heights: header & first data row now match the rest at 13.8
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8

# ---------------------------------------------------------------------------
# 11. Selection ends on G6
# ---------------------------------------------------------------------------
$ws.Range("G6").Select()
